$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing odds cells with new values
$ws.Range("G6").Value = 1.83
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 4.75
$ws.Range("J6").Value = 2.6
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("S6").Value = 1.57
$ws.Range("T6").Value = 2.25
$ws.Range("W6").Value = 5.5
$ws.Range("X6").Value = 7.5
$ws.Range("Z6").Value = 15
$ws.Range("AH6").Value = 9.5
$ws.Range("AI6").Value = 21
$ws.Range("AO6").Value = 11
$ws.Range("AQ6").Value = 41
$ws.Range("AT6").Value = 2.25
$ws.Range("AW6").Value = 6
$ws.Range("AZ6").Value = 101
$ws.Range("G12").Value = 2.35
$ws.Range("N12").Value = 10
$ws.Range("Q12").Value = 2.03
$ws.Range("R12").Value = 1.83
$ws.Range("S12").Value = 1.4
$ws.Range("T12").Value = 2.75
$ws.Range("AT12").Value = 2.75
$ws.Range("BC12").Value = 126
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 11
$ws.Range("Q14").Value = 1.93
$ws.Range("R14").Value = 1.93
$ws.Range("V22").Value = 1.7
$ws.Range("U23").Value = 1.8
$ws.Range("V23").Value = 1.95
$ws.Range("H24").Value = 3.4
$ws.Range("I24").Value = 3
$ws.Range("K24").Value = 2.2
$ws.Range("L24").Value = 3.6
$ws.Range("M24").Value = 1.06
$ws.Range("N24").Value = 10
$ws.Range("O24").Value = 1.29
$ws.Range("P24").Value = 3.5
$ws.Range("Q24").Value = 1.98
$ws.Range("R24").Value = 1.88
$ws.Range("S24").Value = 1.4
$ws.Range("T24").Value = 2.75
$ws.Range("U24").Value = 1.75
$ws.Range("V24").Value = 2
$ws.Range("W24").Value = 8.5
$ws.Range("AB24").Value = 26
$ws.Range("AC24").Value = 10
$ws.Range("AE24").Value = 13
$ws.Range("AF24").Value = 41
$ws.Range("AG24").Value = 201
$ws.Range("AH24").Value = 10
$ws.Range("AK24").Value = 29
$ws.Range("AL24").Value = 23
$ws.Range("AM24").Value = 29
$ws.Range("AP24").Value = 21
$ws.Range("AR24").Value = 51
$ws.Range("AT24").Value = 2.75
$ws.Range("AU24").Value = 7.5
$ws.Range("AY24").Value = 23
$ws.Range("BA24").Value = 67
$ws.Range("BB24").Value = 151
$ws.Range("K25").Value = 2.2
$ws.Range("N25").Value = 13
$ws.Range("O25").Value = 1.25
$ws.Range("P25").Value = 3.75
$ws.Range("Q25").Value = 1.85
$ws.Range("U25").Value = 1.67
$ws.Range("X25").Value = 12
$ws.Range("AM25").Value = 29
$ws.Range("AY25").Value = 23
$ws.Range("J32").Value = 2.38
$ws.Range("L32").Value = 4.33
$ws.Range("U32").Value = 1.67
$ws.Range("V32").Value = 2.1
$ws.Range("W32").Value = 9
$ws.Range("X32").Value = 9.5
$ws.Range("Y32").Value = 8.5
$ws.Range("Z32").Value = 15
$ws.Range("AA32").Value = 13
$ws.Range("AB32").Value = 21
$ws.Range("AD32").Value = 8
$ws.Range("AE32").Value = 15
$ws.Range("AH32").Value = 13
$ws.Range("AI32").Value = 21
$ws.Range("AJ32").Value = 13
$ws.Range("AK32").Value = 41
$ws.Range("AL32").Value = 29
$ws.Range("AM32").Value = 34
$ws.Range("AN32").Value = 4
$ws.Range("AO32").Value = 9
$ws.Range("AP32").Value = 17
$ws.Range("AQ32").Value = 29
$ws.Range("AR32").Value = 41
$ws.Range("AS32").Value = 101
$ws.Range("AW32").Value = 6
$ws.Range("AX32").Value = 21
$ws.Range("AY32").Value = 26
$ws.Range("AZ32").Value = 67
$ws.Range("BA32").Value = 81
$ws.Range("BB32").Value = 151
$ws.Range("M36").Value = 1.02
$ws.Range("O36").Value = 1.13
$ws.Range("M37").Value = 1.03
$ws.Range("O37").Value = 1.17
$ws.Range("M38").Value = 1.03
$ws.Range("O38").Value = 1.17
$ws.Range("M39").Value = 1.03
$ws.Range("O39").Value = 1.17
$ws.Range("G40").Value = 1.62
$ws.Range("I40").Value = 5.75
$ws.Range("J40").Value = 2.25
$ws.Range("L40").Value = 6
$ws.Range("M40").Value = 1.06
$ws.Range("N40").Value = 10
$ws.Range("Q40").Value = 2.07
$ws.Range("R40").Value = 1.69
$ws.Range("U40").Value = 2.1
$ws.Range("V40").Value = 1.67
$ws.Range("AH40").Value = 13
$ws.Range("AI40").Value = 29
$ws.Range("AJ40").Value = 19
$ws.Range("AO40").Value = 8.5
$ws.Range("AS40").Value = 201
$ws.Range("AW40").Value = 7
$ws.Range("AZ40").Value = 126
$ws.Range("G41").Value = 1.6
$ws.Range("H41").Value = 3.8
$ws.Range("I41").Value = 6.25
$ws.Range("M41").Value = 1.05
$ws.Range("N41").Value = 11
$ws.Range("AC41").Value = 11
$ws.Range("AE41").Value = 17
$ws.Range("AO41").Value = 8
$ws.Range("AS41").Value = 151
$ws.Range("Q43").Value = 1.72
$ws.Range("K45").Value = 1.83
$ws.Range("M45").Value = 1.17
$ws.Range("N45").Value = 5
$ws.Range("AC45").Value = 5
$ws.Range("AH45").Value = 10
$ws.Range("AJ45").Value = 23
$ws.Range("AK45").Value = 81
$ws.Range("AL45").Value = 67
$ws.Range("AP45").Value = 34
$ws.Range("BA45").Value = 301
$ws.Range("AH46").Value = 11
$ws.Range("AL46").Value = 21

# Append new row 48 (new match: CF Montreal vs Atlanta Utd)
$ws.Range("A48").Value = "IFBZH1E7"
$ws.Range("B48").Value = "22/10/2024"
$ws.Range("C48").Value = "20:30"
$ws.Range("D48").Value = "USA - MLS"
$ws.Range("E48").Value = "CF Montreal"
$ws.Range("F48").Value = "Atlanta Utd"
$ws.Range("G48").Value = 2.35
$ws.Range("H48").Value = 3.6
$ws.Range("I48").Value = 2.75
$ws.Range("J48").Value = 2.88
$ws.Range("K48").Value = 2.3
$ws.Range("L48").Value = 3.25
$ws.Range("M48").Value = 1.03
$ws.Range("N48").Value = 15
$ws.Range("O48").Value = 1.18
$ws.Range("P48").Value = 4.5
$ws.Range("Q48").Value = 1.62
$ws.Range("R48").Value = 2.25
$ws.Range("S48").Value = 1.3
$ws.Range("T48").Value = 3.4
$ws.Range("U48").Value = 1.53
$ws.Range("V48").Value = 2.38
$ws.Range("W48").Value = 11
$ws.Range("X48").Value = 13
$ws.Range("Y48").Value = 9.5
$ws.Range("Z48").Value = 23
$ws.Range("AA48").Value = 17
$ws.Range("AB48").Value = 21
$ws.Range("AC48").Value = 15
$ws.Range("AD48").Value = 7
$ws.Range("AE48").Value = 12
$ws.Range("AF48").Value = 34
$ws.Range("AG48").Value = 126
$ws.Range("AH48").Value = 12
$ws.Range("AI48").Value = 17
$ws.Range("AJ48").Value = 10
$ws.Range("AK48").Value = 29
$ws.Range("AL48").Value = 21
$ws.Range("AM48").Value = 23
$ws.Range("AN48").Value = 4.75
$ws.Range("AO48").Value = 12
$ws.Range("AP48").Value = 19
$ws.Range("AQ48").Value = 41
$ws.Range("AR48").Value = 51
$ws.Range("AS48").Value = 101
$ws.Range("AT48").Value = 3.4
$ws.Range("AU48").Value = 7
$ws.Range("AV48").Value = 41
$ws.Range("AW48").Value = 5
$ws.Range("AX48").Value = 15
$ws.Range("AY48").Value = 21
$ws.Range("AZ48").Value = 41
$ws.Range("BA48").Value = 51
$ws.Range("BB48").Value = 101
$ws.Range("BC48").Value = 351
$ws.Range("BD48").Value = 151

